$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New training-history values (run with 20 LSTM units, 10 dense units, 50 epochs)
$data = @(
    @(0.09710603207349777, 0.9872792363166809, 0.05815830081701279, 0.9939548373222351),
    @(0.02685911022126675, 0.9978798627853394, 0.04530838504433632, 0.9942730069160461),
    @(0.01106847450137138, 0.9982584714889526, 0.03361079096794128, 0.9955456852912903),
    @(0.007478190120309591, 0.9984666705131531, 0.03296075761318207, 0.9949093461036682),
    @(0.006323980167508125, 0.9985045790672302, 0.07847646623849869, 0.9815462827682495),
    @(0.005134343635290861, 0.9986370801925659, 0.02376365289092064, 0.9958637952804565),
    @(0.004825778771191835, 0.998580276966095, 0.03216518461704254, 0.9958637952804565),
    @(0.004204226657748222, 0.9987884759902954, 0.01572470925748348, 0.9965001344680786),
    @(0.004523877985775471, 0.9986370801925659, 0.01691340282559395, 0.9965001344680786),
    @(0.004420673009008169, 0.9986938238143921, 0.0184360034763813, 0.9971364736557007),
    @(0.004237192682921886, 0.9987127780914307, 0.03463906422257423, 0.9977728128433228),
    @(0.004454357549548149, 0.9987695813179016, 0.02877644076943398, 0.9971364736557007),
    @(0.004041386768221855, 0.998807430267334, 0.02446878328919411, 0.9971364736557007),
    @(0.003855388844385743, 0.9989777803421021, 0.02178826183080673, 0.9974546432495117),
    @(0.004553400445729494, 0.9989399313926697, 0.02268754690885544, 0.9968183040618896),
    @(0.004129610955715179, 0.999034583568573, 0.02228164114058018, 0.9974546432495117),
    @(0.002926370594650507, 0.9994699954986572, 0.02777880802750587, 0.9974546432495117),
    @(0.003564843209460378, 0.9992427825927734, 0.02762977592647076, 0.9974546432495117),
    @(0.004055096302181482, 0.9991102814674377, 0.02798586897552013, 0.9974546432495117),
    @(0.004523102194070816, 0.999034583568573, 0.02222481556236744, 0.9974546432495117),
    @(0.003323583398014307, 0.9993374347686768, 0.02308197878301144, 0.9974546432495117),
    @(0.004548867698758841, 0.9991102814674377, 0.02789357118308544, 0.9984091520309448),
    @(0.004070224240422249, 0.9990913867950439, 0.03223966807126999, 0.9974546432495117),
    @(0.004418177995830774, 0.9991102814674377, 0.02838349901139736, 0.9974546432495117),
    @(0.004044883418828249, 0.9991292357444763, 0.02680773101747036, 0.9971364736557007),
    @(0.004043952561914921, 0.9991670846939087, 0.0321003682911396, 0.9974546432495117),
    @(0.003342764684930444, 0.9993374347686768, 0.03163162618875504, 0.9980909824371338),
    @(0.003697881242260337, 0.999261736869812, 0.02882522717118263, 0.9980909824371338),
    @(0.00396894384175539, 0.9991292357444763, 0.04555697739124298, 0.9977728128433228),
    @(0.003428887110203505, 0.9994131922721863, 0.04260098189115524, 0.9980909824371338),
    @(0.003801741171628237, 0.9991292357444763, 0.03211163356900215, 0.9980909824371338),
    @(0.003890153020620346, 0.9991481900215149, 0.02726985514163971, 0.9980909824371338),
    @(0.003542009741067886, 0.999261736869812, 0.03183046728372574, 0.9980909824371338),
    @(0.004088382236659527, 0.9991102814674377, 0.0330834835767746, 0.9980909824371338),
    @(0.003787299152463675, 0.9992238879203796, 0.03787132725119591, 0.9974546432495117),
    @(0.00311693549156189, 0.9994131922721863, 0.03463936224579811, 0.9974546432495117),
    @(0.003718815045431256, 0.9991860389709473, 0.03611919283866882, 0.9977728128433228),
    @(0.003628178266808391, 0.9992427825927734, 0.05458014458417892, 0.9974546432495117),
    @(0.003968276549130678, 0.9990913867950439, 0.02719040773808956, 0.9980909824371338),
    @(0.005423769820481539, 0.999034583568573, 0.05385555699467659, 0.9977728128433228),
    @(0.003552494570612907, 0.999318540096283, 0.05979811027646065, 0.9974546432495117),
    @(0.003552684094756842, 0.999261736869812, 0.05861984938383102, 0.9974546432495117),
    @(0.00392773887142539, 0.9992049336433411, 0.04595331102609634, 0.9980909824371338),
    @(0.003732847049832344, 0.999261736869812, 0.03814717754721642, 0.9974546432495117),
    @(0.002643566112965345, 0.9995078444480896, 0.04010450839996338, 0.9974546432495117),
    @(0.002710732631385326, 0.9994510412216187, 0.04469374567270279, 0.9980909824371338),
    @(0.003474482102319598, 0.999261736869812, 0.04085085168480873, 0.9980909824371338),
    @(0.003149901051074266, 0.9993563890457153, 0.04933448880910873, 0.9974546432495117),
    @(0.004067272879183292, 0.9992049336433411, 0.04717199131846428, 0.9974546432495117),
    @(0.004523102194070816, 0.9992806911468506, 0.04833652451634407, 0.9971364736557007)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
